$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price + Volume(1h) refresh from the data
# source), including a few rows where the coin at that rank changed so
# name/link/price/volume were all replaced.
#
# Cells whose new text would otherwise be auto-recognised by Excel as a
# number (plain decimals like "577.66") are forced back to Text first so
# they round-trip as the exact literal string from the source diff,
# matching the original inlineStr cells (which are all text, never real
# numbers, in this sheet).

$ws.Range('D2').Value = '67.291.75'
$ws.Range('E2').Value = '  -1.39%  '

$ws.Range('D3').Value = '3.589.95'
$ws.Range('E3').Value = '  -2.70%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.66'
$ws.Range('E5').Value = '  -4.98%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '192.33'
$ws.Range('E6').Value = '  -0.31%  '

$ws.Range('D7').Value = '3.586.99'
$ws.Range('E7').Value = '  -2.65%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.618'
$ws.Range('E8').Value = '  -2.35%  '

$ws.Range('E9').Value = '  +0.15%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.680'
$ws.Range('E10').Value = '  -5.90%  '

$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '56.18'
$ws.Range('E11').Value = '  -5.35%  '

$ws.Range('B12').Value = 'Dogecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.150'
$ws.Range('E12').Value = '  -5.56%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000274'
$ws.Range('E13').Value = '  -4.09%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.86'
$ws.Range('E14').Value = '  -5.17%  '

$ws.Range('D15').Value = '4.175.08'
$ws.Range('E15').Value = '  -2.61%  '

$ws.Range('D16').Value = '3.597.42'
$ws.Range('E16').Value = '  -2.72%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.41'
$ws.Range('E18').Value = '  -5.04%  '

$ws.Range('D19').Value = '67.307.31'
$ws.Range('E19').Value = '  -1.29%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.19'
$ws.Range('E20').Value = '  -4.73%  '

$ws.Range('E21').Value = '  -6.65%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '402.65'
$ws.Range('E22').Value = '  -1.09%  '

$ws.Range('E23').Value = '  -7.71%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.93'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.39'
$ws.Range('E25').Value = '  -0.85%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.94'
$ws.Range('E26').Value = '  -4.11%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.49'
$ws.Range('E27').Value = '  -3.78%  '

$ws.Range('E28').Value = '  +1.15%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.65'
$ws.Range('E29').Value = '  -2.74%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.97'
$ws.Range('E30').Value = '  -6.40%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.65'
$ws.Range('E31').Value = '  -0.30%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '31.25'
$ws.Range('E32').Value = '  -3.99%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '635.18'
$ws.Range('E33').Value = '  +0.71%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '12.20'
$ws.Range('E34').Value = '  -3.78%  '

$ws.Range('E35').Value = '  -5.34%  '

$ws.Range('E36').Value = '  -4.75%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '42.49'
$ws.Range('E37').Value = '  -10.70%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.399'
$ws.Range('E38').Value = '  -2.94%  '

$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +0.21%  '

$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D40').Value = '0.0₃0785'
$ws.Range('E40').Value = '  -4.39%  '

$ws.Range('D41').Value = '3.173.84'
$ws.Range('E41').Value = '  +8.57%  '

$ws.Range('E42').Value = '  -2.99%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.73'
$ws.Range('E43').Value = '  +4.36%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').Value = '  -0.24%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.00'
$ws.Range('E45').Value = '  -0.90%  '

$ws.Range('E46').Value = '  -5.60%  '

$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.12'
$ws.Range('E47').Value = '  +0.41%  '

$ws.Range('E48').Value = '  -6.23%  '

$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.64'
$ws.Range('E49').Value = '  -2.50%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '141.73'
$ws.Range('E50').Value = '  -3.02%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.61'
$ws.Range('E51').Value = '  -6.83%  '
